# Apply "gh-pages output generated at 456a3b4" update to 杭州-漫展信息.xlsx
#
# Summary of changes:
#  - Sheet "展览"   (Exhibition):  bump several "想去人数" (want-to-go) F counts,
#                                   and fix G18 lowest-price from 59.9 to 75.
#  - Sheet "演出"   (Performance): bump F3 want-to-go count.
#  - Sheet "本地生活" (Local Life): the oldest listing (剑网3×HAPPY ZOO) expired and
#                                   was removed; all rows below shift up one, and
#                                   the want-to-go count for the newest listing
#                                   (蜡笔小新) ticks up from 143 to 144.
#  - Sheet "全部类型" (All Types):  mirrors the same F/G bumps as 展览/演出 (this
#                                   sheet does not include 本地生活 rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2677
$ws1.Range("F5").Value = 927
$ws1.Range("F7").Value = 2177
$ws1.Range("F8").Value = 1802
$ws1.Range("F9").Value = 209
$ws1.Range("F11").Value = 2458
$ws1.Range("F13").Value = 228
$ws1.Range("F16").Value = 122
$ws1.Range("F18").Value = 9108
$ws1.Range("G18").Value = 75
$ws1.Range("F20").Value = 7058
$ws1.Range("F21").Value = 11546
$ws1.Range("F24").Value = 232
$ws1.Range("F25").Value = 339
$ws1.Range("F26").Value = 549
$ws1.Range("F27").Value = 2545
$ws1.Range("F28").Value = 229
$ws1.Range("F29").Value = 194
$ws1.Range("F30").Value = 2481
$ws1.Range("F31").Value = 653
$ws1.Range("F32").Value = 44
$ws1.Range("F33").Value = 4495
$ws1.Range("F34").Value = 843
$ws1.Range("F35").Value = 343
$ws1.Range("F36").Value = 38
$ws1.Range("F37").Value = 511

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 68

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life)
#   Row 2 (2024-08-27 剑网3xHAPPY ZOO) is removed; everything below shifts up.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Rows.Item(2).Delete()

# Re-number the index column (A) back to the sequential 1,2,3 it had before
# the delete shifted the old row 3/4/5 index values up into rows 2/3/4.
$ws3.Range("A2").Value = 1
$ws3.Range("A3").Value = 2
$ws3.Range("A4").Value = 3

# The want-to-go count for the (now last) 蜡笔小新 listing ticked up.
$ws3.Range("F4").Value = 144

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2677
$ws4.Range("F6").Value = 68
$ws4.Range("F7").Value = 927
$ws4.Range("F9").Value = 2177
$ws4.Range("F11").Value = 1802
$ws4.Range("F13").Value = 209
$ws4.Range("F14").Value = 2458
$ws4.Range("F17").Value = 228
$ws4.Range("F20").Value = 122
$ws4.Range("F22").Value = 9108
$ws4.Range("G22").Value = 75
$ws4.Range("F24").Value = 7058
$ws4.Range("F25").Value = 11546
$ws4.Range("F28").Value = 232
$ws4.Range("F29").Value = 339
$ws4.Range("F31").Value = 549
$ws4.Range("F33").Value = 2545
$ws4.Range("F36").Value = 229
$ws4.Range("F37").Value = 194
$ws4.Range("F38").Value = 44
$ws4.Range("F39").Value = 4495
$ws4.Range("F46").Value = 511
